# sprint 1 & 2
# Updates the "Sprint 3" grading table on the "Group 1" worksheet:
#  - removes the "Backlog has user story 5 tasks" criterion row
#  - removes the "User story 5 implemented" criterion row
#  - inserts a new "Test database configuration" criterion row
#    (right after "README has license and test usage guide")
#  - bumps the max points (column D) for the "User story 1-4 implemented"
#    rows from 0.75 to 1
#  - reduces the max points for "User stories from previous Sprints
#    implemented" from 1 to 0.75
#  - refreshes the view (zoom / scroll / selection) on that sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group 1")

# --- Remove "Backlog has user story 5 tasks" (row 63) ---------------------
$ws.Rows.Item(63).Delete()

# --- Remove "User story 5 implemented" (now row 71 after the shift) -------
$ws.Rows.Item(71).Delete()

# --- Insert the new "Test database configuration" row after row 65 --------
# (row 65 is now "README has license and test usage guide")
$ws.Rows.Item(66).Insert()
$ws.Range("A66").Value = "Test database configuration"
$ws.Range("C66").Value = 0
$ws.Range("D66").Value = 0.25

# --- Update max points for the user story implementation rows -------------
$ws.Range("D68").Value = 1
$ws.Range("D69").Value = 1
$ws.Range("D70").Value = 1
$ws.Range("D71").Value = 1

# --- Update max points for "User stories from previous Sprints implemented"
$ws.Range("D76").Value = 0.75

# --- Refresh the sheet view (scroll position / zoom / selection) ----------
[void]$ws.Activate()
[void]$ws.Range("A66").Select()
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 95

Write-Host "Sprint 3 criteria updated."
